# Query Log workbook update — "update model & more"
#
# Adds a small H/M/L priority "model" list (rows 25-28) used to drive a new
# list-based data validation on the Query Priority column (C9:C16), fills in
# the priority + resolution-date values that were picked from that list,
# adds a new query row (row 16) and updates the current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Query Log")
$ws.Activate()

# --- New "model" lookup table (Setup-style list), rows 25-28 ---------------
# Written H (row26) before status (row25) so new shared-string indices land
# in the same order as the authored workbook: H, status, M, L.
$ws.Range("A26").Value = "H"
$ws.Range("A25").Value = "status"
$ws.Range("A27").Value = "M"
$ws.Range("A28").Value = "L"

# --- List data validation on C9:C16, sourced from the new A26:A28 list -----
$priorityRange = $ws.Range("C9:C16")
$priorityRange.Validation.Delete()
$priorityRange.Validation.Add(3, 1, 1, "=A26:A28")
$priorityRange.Validation.IgnoreBlank = $true
$priorityRange.Validation.InCellDropdown = $true
$priorityRange.Validation.ShowInput = $true
$priorityRange.Validation.ShowError = $true

# --- Fill in Query Priority selections + Resolution Date for existing rows -
$ws.Range("C9").Value = "H"

$ws.Range("C10").Value = "L"
$ws.Range("J10").Value = "9/14/2020"

$ws.Range("C11").Value = "L"
$ws.Range("J11").Value = "9/14/2020"

$ws.Range("C12").Value = "L"
$ws.Range("J12").Value = "9/14/2020"

$ws.Range("C13").Value = "L"
$ws.Range("J13").Value = "9/14/2020"

$ws.Range("C14").Value = "L"
$ws.Range("J14").Value = "9/14/2020"

$ws.Range("C15").Value = "H"

# --- New query row 16 -------------------------------------------------------
$ws.Rows.Item(16).RowHeight = 154

$ws.Range("B16").Value = "should the user be able to place a order from the client side or orders are only to come from external source e.g manually inserting into the database ?"
$ws.Range("C16").Value = "H"
$ws.Range("D16").Value = "Buchi"
$ws.Range("E16").Value = "9/14/2020"
$ws.Range("F16").Value = "instructers "
$ws.Range("G16").Value = "ASAP"

# --- Selection matches the author's final cursor position ------------------
$ws.Range("I15").Select()
